# "move aws iam into aws iam labs"
#
# The "Hands On: IAM" slide (SlideID 258 - the last slide in the deck,
# along with its notes page) is being moved out of this deck and into a
# separate "aws iam labs" deck, so here it is simply removed.

$p = $ppt.ActivePresentation

# Find the slide with SlideID 258 ("Hands On: IAM") rather than assuming
# a fixed index, then delete it. Deleting the slide also removes its
# associated notes page (notesSlideN.xml) and its entry in the slide
# id list.
$targetId = 258
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq $targetId) {
        $slide.Delete()
        break
    }
}
